$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue "D2" "28.918.31"
Set-TextValue "E2" "  +1.68%  "
Set-TextValue "D3" "1.890.44"
Set-TextValue "E3" "  +1.58%  "
Set-TextValue "D4" "1.001"
Set-TextValue "E4" "  -0.67%  "
Set-TextValue "D5" "325.24"
Set-TextValue "E5" "  -0.08%  "
Set-TextValue "D6" "1.000"
Set-TextValue "E6" "  -0.66%  "
Set-TextValue "D7" "0.4582"
Set-TextValue "E7" "  +0.58%  "
Set-TextValue "D8" "0.3904"
Set-TextValue "E8" "  +2.02%  "
Set-TextValue "D9" "0.07837"
Set-TextValue "E9" "  +0.36%  "
Set-TextValue "D10" "0.9881"
Set-TextValue "E10" "  +0.37%  "
Set-TextValue "D11" "21.87"
Set-TextValue "E11" "  +1.91%  "
Set-TextValue "D12" "1.942.96"
Set-TextValue "E12" "  +5.39%  "
Set-TextValue "D13" "7.048"
Set-TextValue "E13" "  +2.33%  "
Set-TextValue "D14" "5.685"
Set-TextValue "E14" "  +1.00%  "
Set-TextValue "D15" "0.06929"
Set-TextValue "E15" "  +0.29%  "
Set-TextValue "D16" "87.95"
Set-TextValue "E16" "  +1.82%  "
Set-TextValue "D17" "1.002"
Set-TextValue "E17" "  -0.65%  "
Set-TextValue "D18" "0.000009972"
Set-TextValue "E18" "  +0.47%  "
Set-TextValue "D19" "16.98"
Set-TextValue "E19" "  +2.00%  "
Set-TextValue "D20" "1.002"
Set-TextValue "E20" "  -0.44%  "
Set-TextValue "D21" "28.944.76"
Set-TextValue "E21" "  +1.77%  "
Set-TextValue "D22" "5.291"
Set-TextValue "E22" "  +0.87%  "
Set-TextValue "D23" "10.98"
Set-TextValue "E23" "  +1.08%  "
Set-TextValue "D24" "2.161.65"
Set-TextValue "E24" "  +4.49%  "
Set-TextValue "D25" "2.060"
Set-TextValue "E25" "  -1.51%  "
Set-TextValue "D26" "156.10"
Set-TextValue "E26" "  +2.01%  "
Set-TextValue "D27" "19.25"
Set-TextValue "E27" "  +0.88%  "
Set-TextValue "D28" "5.911"
Set-TextValue "E28" "  +4.81%  "
Set-TextValue "D29" "1.928"
Set-TextValue "E29" "  +1.67%  "
Set-TextValue "D30" "117.43"
Set-TextValue "E30" "  +0.29%  "
Set-TextValue "D31" "0.09341"
Set-TextValue "D32" "0.9053"
Set-TextValue "E32" "  +0.25%  "
Set-TextValue "D33" "5.291"
Set-TextValue "E33" "  +0.50%  "
Set-TextValue "D34" "1.330"
Set-TextValue "E34" "  +1.24%  "
Set-TextValue "D35" "3.258"
Set-TextValue "E35" "  -1.00%  "
Set-TextValue "D36" "1.190"
Set-TextValue "E36" "  +3.30%  "
Set-TextValue "D37" "0.05767"
Set-TextValue "E37" "  +1.56%  "
Set-TextValue "D38" "0.02074"
Set-TextValue "E38" "  +1.97%  "
Set-TextValue "D39" "0.9998"
Set-TextValue "E39" "  -0.64%  "
Set-TextValue "D40" "7.742"
Set-TextValue "E40" "  +1.42%  "
Set-TextValue "D41" "0.5680"
Set-TextValue "E41" "  +2.56%  "
Set-TextValue "D42" "0.1771"
Set-TextValue "E42" "  +0.47%  "
Set-TextValue "D43" "9.733"
Set-TextValue "E43" "  +1.46%  "
Set-TextValue "D44" "2.298"
Set-TextValue "D45" "11.97"
Set-TextValue "E45" "  +4.29%  "
Set-TextValue "D46" "0.5343"
Set-TextValue "E46" "  +2.29%  "
Set-TextValue "E47" "  -1.10%  "
Set-TextValue "D48" "1.845"
Set-TextValue "E48" "  +2.37%  "
Set-TextValue "D49" "112.80"
Set-TextValue "E49" "  +1.04%  "
Set-TextValue "D50" "2.527"
Set-TextValue "E50" "  +3.78%  "
Set-TextValue "E51" "  -5.32%  "
